$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.03864900479865655
$ws.Range("C2").Value = 0.5850438662025559
$ws.Range("D2").Value = 0.7245914535781601
$ws.Range("E2").Value = 0.8512293777696821
$ws.Range("F2").Value = 0.8736531364880852
$ws.Range("G2").Value = 19
$ws.Range("B3").Value = -0.05426294671939138
$ws.Range("C3").Value = 0.5223775699184166
$ws.Range("D3").Value = 0.5969843308118188
$ws.Range("E3").Value = 0.7726476110179975
$ws.Range("F3").Value = 0.7930847273798034
$ws.Range("G3").Value = 18
$ws.Range("B4").Value = 0.04121462200618248
$ws.Range("C4").Value = 0.4282064007022935
$ws.Range("D4").Value = 0.4114353831211015
$ws.Range("E4").Value = 0.641432290363606
$ws.Range("F4").Value = 0.6598070052540843
$ws.Range("G4").Value = 17
$ws.Range("B5").Value = 0.1513232839577969
$ws.Range("C5").Value = 0.3552284090837819
$ws.Range("D5").Value = 0.2325151131945527
$ws.Range("E5").Value = 0.4821982094476842
$ws.Range("F5").Value = 0.4728538908111392
$ws.Range("G5").Value = 16
$ws.Range("B6").Value = 0.1388059945772336
$ws.Range("C6").Value = 0.3804343702871964
$ws.Range("D6").Value = 0.1997170828207708
$ws.Range("E6").Value = 0.4468971725361113
$ws.Range("F6").Value = 0.4397036079933307
$ws.Range("G6").Value = 15
$ws.Range("B7").Value = 0.1640011163335129
$ws.Range("C7").Value = 0.3317361490081165
$ws.Range("D7").Value = 0.2137299874398079
$ws.Range("E7").Value = 0.4623094066096946
$ws.Range("F7").Value = 0.4485592918475749
$ws.Range("G7").Value = 14
$ws.Range("B8").Value = 0.1805696058411114
$ws.Range("C8").Value = 0.3326462541290968
$ws.Range("D8").Value = 0.1927252227720784
$ws.Range("E8").Value = 0.4390048095090513
$ws.Range("F8").Value = 0.4164890877762218
$ws.Range("G8").Value = 13
$ws.Range("B9").Value = 0.217823495412195
$ws.Range("C9").Value = 0.3128752448572103
$ws.Range("D9").Value = 0.2168929369963322
$ws.Range("E9").Value = 0.4657176580250443
$ws.Range("F9").Value = 0.4299418927031617
$ws.Range("G9").Value = 12
$ws.Range("B10").Value = 0.1640088941767343
$ws.Range("C10").Value = 0.2728317493205706
$ws.Range("D10").Value = 0.1256684404680174
$ws.Range("E10").Value = 0.3544974477595254
$ws.Range("F10").Value = 0.3296156480036049
$ws.Range("G10").Value = 11
$ws.Range("B11").Value = 0.1230033177014676
$ws.Range("C11").Value = 0.3182986834576886
$ws.Range("D11").Value = 0.2493480347795076
$ws.Range("E11").Value = 0.499347609165707
$ws.Range("F11").Value = 0.5101396525723122
$ws.Range("G11").Value = 10
